# NSMB - 5-2 done. Extend the "V3" run-comparison sheet (sheet1 / first
# tab) with new checkpoint rows 49-59, tweak row 48 (saved 48 frames: final
# clock time improved from 19150 to 19152, so diffs drop 242->240/18->16),
# and move the frozen-pane selection down to B60.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 48 update: one-off B value correction and recomputed diffs ---
$ws.Range("B48").Value = 19152

# --- Row 49 ---
$ws.Range("A49").Value = "Checkpoint 732"
$ws.Range("B49").Value = 19206
$ws.Range("C49").Value = 19447
$ws.Range("E49").Value = 19223

# --- Row 50 ---
$ws.Range("A50").Value = "Checkpoint 949"
$ws.Range("B50").Value = 19281
$ws.Range("C50").Value = 19522
$ws.Range("E50").Value = 19299

# --- Row 51 ---
$ws.Range("A51").Value = "Checkpoint 1154"
$ws.Range("B51").Value = 19350
$ws.Range("C51").Value = 19591

# --- Row 52 ---
$ws.Range("A52").Value = "Checkpoint 1221"
$ws.Range("B52").Value = 19374
$ws.Range("C52").Value = 19615

# --- Row 53 ---
$ws.Range("A53").Value = "Checkpoint 1470"
$ws.Range("B53").Value = 19461
$ws.Range("C53").Value = 19703

# --- Row 54 ---
$ws.Range("A54").Value = "Checkpoint 2208"
$ws.Range("B54").Value = 19706
$ws.Range("C54").Value = 19949

# --- Row 55 ---
$ws.Range("A55").Value = "Checkpoitn 3427 (1st time)"
$ws.Range("B55").Value = 20114
$ws.Range("C55").Value = 20359

# --- Row 56 ---
$ws.Range("A56").Value = "Checkpoint 3422 (2nd time)"
$ws.Range("B56").Value = 20117
$ws.Range("C56").Value = 20363

# --- Row 57 (reuses existing "Enter pipe" shared string) ---
$ws.Range("A57").Value = "Enter pipe"
$ws.Range("B57").Value = 20263
$ws.Range("C57").Value = 20512

# --- Row 58 (reuses existing "Get flag" shared string) ---
$ws.Range("A58").Value = "Get flag"
$ws.Range("B58").Value = 20543
$ws.Range("C58").Value = 20832

# --- Row 59 (reuses existing "Black screen" shared string) ---
$ws.Range("A59").Value = "Black screen"
$ws.Range("B59").Value = 21058
$ws.Range("C59").Value = 21346

# --- Fill the D (C-B) shared formula down through the new rows ---
$ws.Range("D49:D59").Formula = "=IF(B49 >  0,C49-B49, 0)"

# --- Fill the F (E-B) shared formula down through the two new rows that
#     still have an E value (49 and 50), matching the source ref F39:F50 ---
$ws.Range("F49:F50").Formula = "=IF(B49 >  0,E49-B49, 0)"

# --- Recalculate so cached <v> results land in the saved workbook ---
$excel.Calculate()

# --- Move the selection to follow the newly-added last row (B60), as the
#     author's saved view shows ---
$ws.Range("B60").Select()
